$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue 'D2' '29.392.02'
Set-TextValue 'E2' '  +0.28%  '
Set-TextValue 'D3' '1.842.93'
Set-TextValue 'E3' '  +0.08%  '
Set-TextValue 'E4' '  -0.04%  '
Set-TextValue 'D5' '240.09'
Set-TextValue 'E5' '  -0.06%  '
Set-TextValue 'D6' '0.6336'
Set-TextValue 'E6' '  +1.04%  '
Set-TextValue 'D7' '0.9998'
Set-TextValue 'E7' '  -0.10%  '
Set-TextValue 'D8' '0.07473'
Set-TextValue 'E8' '  -0.13%  '
Set-TextValue 'B9' 'Solana'
Set-TextValue 'C9' 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
Set-TextValue 'D9' '25.08'
Set-TextValue 'E9' '  +3.12%  '
Set-TextValue 'B10' 'Cardano'
Set-TextValue 'C10' 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
Set-TextValue 'D10' '0.2905'
Set-TextValue 'E10' '  +0.36%  '
Set-TextValue 'E11' '  +0.42%  '
Set-TextValue 'D12' '1.849.02'
Set-TextValue 'E12' '  +0.38%  '
Set-TextValue 'D13' '4.989'
Set-TextValue 'E13' '  +0.05%  '
Set-TextValue 'D14' '0.6792'
Set-TextValue 'E14' '  +0.24%  '
Set-TextValue 'D15' '0.00001022'
Set-TextValue 'E15' '  -0.32%  '
Set-TextValue 'D16' '82.02'
Set-TextValue 'E16' '  -0.09%  '
Set-TextValue 'D17' '6.268'
Set-TextValue 'E17' '  +2.69%  '
Set-TextValue 'D18' '29.423.05'
Set-TextValue 'E18' '  +0.27%  '
Set-TextValue 'D19' '230.26'
Set-TextValue 'E19' '  +0.80%  '
Set-TextValue 'D20' '12.34'
Set-TextValue 'E20' '  +0.71%  '
Set-TextValue 'D21' '0.9996'
Set-TextValue 'E21' '  -0.09%  '
Set-TextValue 'D22' '7.424'
Set-TextValue 'E22' '  +0.62%  '
Set-TextValue 'D23' '0.9999'
Set-TextValue 'E23' '  -0.05%  '
Set-TextValue 'D24' '158.17'
Set-TextValue 'E24' '  -0.32%  '
Set-TextValue 'D25' '8.502'
Set-TextValue 'E25' '  +1.54%  '
Set-TextValue 'D26' '0.1358'
Set-TextValue 'E26' '  -1.66%  '
Set-TextValue 'D27' '17.48'
Set-TextValue 'E27' '  -0.33%  '
Set-TextValue 'D28' '0.06546'
Set-TextValue 'E28' '  +14.76%  '
Set-TextValue 'D29' '1.431'
Set-TextValue 'E29' '  +2.57%  '
Set-TextValue 'D30' '1.489'
Set-TextValue 'E30' '  +1.02%  '
Set-TextValue 'E31' '  -0.48%  '
Set-TextValue 'E32' '  +0.51%  '
Set-TextValue 'E33' '  +1.37%  '
Set-TextValue 'D34' '1.141'
Set-TextValue 'E34' '  -0.09%  '
Set-TextValue 'D35' '0.6982'
Set-TextValue 'E35' '  +1.10%  '
Set-TextValue 'E37' '  +2.53%  '
Set-TextValue 'E38' '  -0.65%  '
Set-TextValue 'D39' '1.248.04'
Set-TextValue 'E39' '  +0.15%  '
Set-TextValue 'D40' '6.771'
Set-TextValue 'E40' '  +3.93%  '
Set-TextValue 'E41' '  +3.23%  '
Set-TextValue 'D42' '0.9997'
Set-TextValue 'E42' '  +0.09%  '
Set-TextValue 'D43' '2.006.30'
Set-TextValue 'E43' '  +0.21%  '
Set-TextValue 'D44' '101.14'
Set-TextValue 'E44' '  -0.16%  '
Set-TextValue 'D45' '65.41'
Set-TextValue 'E45' '  -0.42%  '
Set-TextValue 'B46' 'BabyDogeCoin'
Set-TextValue 'C46' 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextValue 'D46' '0.00000000120'
Set-TextValue 'E46' '  +3.40%  '
Set-TextValue 'B47' 'Aptos'
Set-TextValue 'C47' 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue 'D47' '7.066'
Set-TextValue 'E47' '  -0.11%  '
Set-TextValue 'D48' '1.715'
Set-TextValue 'E48' '  +3.70%  '
Set-TextValue 'D49' '9.023'
Set-TextValue 'E49' '  +0.18%  '
Set-TextValue 'D50' '0.1147'
Set-TextValue 'E50' '  -1.26%  '
Set-TextValue 'E51' '  -0.54%  '
